# Update "Coding approach" slide (slide 6) coding-approach bullet list:
# add "Leaflet", "Plotly", ".GeoJSON" paragraphs after "Flask app "
$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(6)
$contentShape = $slide.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange

# Paragraph 4 is "Flask app " - insert the new paragraphs right after it,
# one at a time, always re-fetching the freshly created last paragraph so
# InsertAfter appends in order instead of stacking in reverse.
$flaskPara = $tr.Paragraphs(4, 1)
$flaskPara.InsertAfter("`rLeaflet") | Out-Null

$leafletPara = $tr.Paragraphs(5, 1)
$leafletPara.InsertAfter("`rPlotly") | Out-Null

$plotlyPara = $tr.Paragraphs(6, 1)
$plotlyPara.InsertAfter("`r.") | Out-Null

$dotPara = $tr.Paragraphs(7, 1)
$dotPara.InsertAfter("GeoJSON") | Out-Null

# Update every cached "datetimeFigureOut" date field (slide master, notes
# master, and all slide layouts) from the long format (8/27/2020) to the
# short format (8/27/20).
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $dateRange = $sh.TextFrame.TextRange
            if ($dateRange.Text -eq "8/27/2020") {
                $dateRange.Text = "8/27/20"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes
Update-DatePlaceholder $p.NotesMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
